$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Select header row on sheet1 (matches the target's sheetView selection)
$ws1.Range("A1:E1").Select()

# Add the new sheet "16 bat database matches" right after the existing sheet
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "16 bat database matches"

# Header row, same layout/labels as sheet1
$ws2.Range("A1").Value = "accession number"
$ws2.Range("B1").Value = "human db p-value"
$ws2.Range("C1").Value = "human db mean query PCC"
$ws2.Range("D1").Value = "bat db p-value"
$ws2.Range("E1").Value = "bat db mean PCC"
$ws2.Range("A1:E1").Font.Bold = $true
$ws2.Range("A1:E1").Interior.Color = 49407

$ws2.Range("A2").Value = "GQ153546.1"
$ws2.Range("B2").Value = "4.0663387686233697E-241"
$ws2.Range("C2").Value = "0.65696195322147299"
$ws2.Range("D2").Value = "0.033851270588193198"
$ws2.Range("E2").Value = "0.85293975471199801"

$ws2.Range("A3").Value = "GQ153547.1"
$ws2.Range("B3").Value = "9.7111592768242005E-248"
$ws2.Range("C3").Value = "0.66036964706026102"
$ws2.Range("D3").Value = "0.061142171505474401"
$ws2.Range("E3").Value = "0.84404136924456197"

$ws2.Range("A4").Value = "GQ153548.1"
$ws2.Range("B4").Value = "1.2080639484103599E-239"
$ws2.Range("C4").Value = "0.65400993315080602"
$ws2.Range("D4").Value = "0.039417814644244001"
$ws2.Range("E4").Value = "0.85180692404284897"

$ws2.Range("A5").Value = "GU190215.1"
$ws2.Range("B5").Value = "0"
$ws2.Range("C5").Value = "0.63815838441720596"
$ws2.Range("D5").Value = "7.6231893973146102E-15"
$ws2.Range("E5").Value = "0.672520775875874"

$ws2.Range("A6").Value = "JX993987.1"
$ws2.Range("B6").Value = "6.7312026862475896E-224"
$ws2.Range("C6").Value = "0.67018595418517402"
$ws2.Range("D6").Value = "0.0096014101720237992"
$ws2.Range("E6").Value = "0.78368965149940095"

$ws2.Range("A7").Value = "JX993988.1"
$ws2.Range("B7").Value = "4.2287847510923799E-193"
$ws2.Range("C7").Value = "0.67901372102996804"
$ws2.Range("D7").Value = "0.0026149079601012599"
$ws2.Range("E7").Value = "0.77626189129910805"

$ws2.Range("A8").Value = "KC881005.1"
$ws2.Range("B8").Value = "7.0334241064054105E-128"
$ws2.Range("C8").Value = "0.67235521368317996"
$ws2.Range("D8").Value = "0.82150903555737198"
$ws2.Range("E8").Value = "0.811349465999609"

$ws2.Range("A9").Value = "KC881006.1"
$ws2.Range("B9").Value = "1.07861556109611E-125"
$ws2.Range("C9").Value = "0.67384871426196202"
$ws2.Range("D9").Value = "0.69906642808255504"
$ws2.Range("E9").Value = "0.813728775403041"

$ws2.Range("A10").Value = "KF294457.1"
$ws2.Range("B10").Value = "1.7001213861528199E-306"
$ws2.Range("C10").Value = "0.66190982472336002"
$ws2.Range("D10").Value = "0.75452557209323201"
$ws2.Range("E10").Value = "0.80416129805834102"

$ws2.Range("A11").Value = "KF367457.1"
$ws2.Range("B11").Value = "1.9199869277680499E-127"
$ws2.Range("C11").Value = "0.67149897645283396"
$ws2.Range("D11").Value = "0.83407295214907096"
$ws2.Range("E11").Value = "0.81109262990992403"

$ws2.Range("A12").Value = "KP886808.1"
$ws2.Range("B12").Value = "3.7967318383287897E-141"
$ws2.Range("C12").Value = "0.67403404295014402"
$ws2.Range("D12").Value = "0.38972362870238297"
$ws2.Range("E12").Value = "0.79691611292239795"

$ws2.Range("A13").Value = "KP886809.1"
$ws2.Range("B13").Value = "1.43587096264173E-141"
$ws2.Range("C13").Value = "0.67205103113254805"
$ws2.Range("D13").Value = "0.33991191476485899"
$ws2.Range("E13").Value = "0.79561467413863596"

$ws2.Range("A14").Value = "KU182964.1"
$ws2.Range("B14").Value = "0"
$ws2.Range("B14").NumberFormat = "0.00E+00"
$ws2.Range("C14").Value = "0.66708791010942603"
$ws2.Range("D14").Value = "0.00102563400789871"
$ws2.Range("E14").Value = "0.76837804933429599"

$ws2.Range("A15").Value = "KY417142.1"
$ws2.Range("B15").Value = "3.60696420883077E-133"
$ws2.Range("C15").Value = "0.67179533390970503"
$ws2.Range("D15").Value = "0.86421118411890996"
$ws2.Range("E15").Value = "0.81034989445002603"

$ws2.Range("A16").Value = "KY417143.1"
$ws2.Range("B16").Value = "8.2754196174979096E-136"
$ws2.Range("C16").Value = "0.67872856338883003"
$ws2.Range("D16").Value = "0.84612269898944203"
$ws2.Range("E16").Value = "0.810592559109795"

$ws2.Range("A17").Value = "KY417144.1"
$ws2.Range("B17").Value = "9.7540644659330492E-127"
$ws2.Range("C17").Value = "0.67128177285591095"
$ws2.Range("D17").Value = "0.81276020569865404"
$ws2.Range("E17").Value = "0.81143181474509096"


# Final selection on the new sheet (matches target's sheetView selection on sheet2)
$ws2.Range("C2").Select()
